$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 to shift existing data down by one row.
$ws.Rows.Item(2).Insert()

# The insert copies formatting (including the wrap-text style) from the row
# above; clear it so the new row matches the plain, unstyled data rows.
$ws.Rows.Item(2).ClearFormats()

# Populate the newly inserted row 2 with the "cruise" attribute metadata.
$ws.Range("A2").Value = "cruise"
$ws.Range("B2").Value = "Identifier for research cruise generally including abbreviation for research vessel and voyage number"
$ws.Range("C2").Value = "character"

# Update datetime_utc_matlab row (now row 4): class should be Date with a
# dateTimeFormatString instead of numeric/dimensionless.
$ws.Range("C4").Value = "Date"
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = "YYYY-MM-DD hh:mm:ss"

# Update the active selection to match the saved view state.
$ws.Range("E7").Select()
